# Auto-generated edit script: updates crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''67.259.38'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  +5.81%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''3.719.90'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  +7.22%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = '''  +0.04%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''424.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  +2.32%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''131.94'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  +2.18%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = '''3.712.03'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''  +7.28%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = '''0.644'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''  +2.93%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = '''  +0.00%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''0.770'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  -1.45%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''0.184'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  +13.06%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''0.0000394'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  +56.58%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''43.10'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  +1.60%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''10.20'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  +3.80%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''4.295.77'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  +6.88%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('E16').Value = '''  -0.01%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = '''20.87'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  +3.21%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''3.765.91'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  +8.60%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('B19').Value = '''Polygon'
$ws.Range('B19').Style = 'Normal'
$ws.Range('C19').Value = '''https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('C19').Style = 'Normal'
$ws.Range('D19').Value = '''1.14'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  +4.87%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('B20').Value = '''Uniswap'
$ws.Range('B20').Style = 'Normal'
$ws.Range('C20').Value = '''https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('C20').Style = 'Normal'
$ws.Range('D20').Value = '''12.91'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  +3.82%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''67.265.76'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  +6.07%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''451.47'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  -2.47%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''15.95'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  +18.00%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''89.83'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  -0.39%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''3.21'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  -2.56%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = '''38.08'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  +12.20%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = '''10.28'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  +1.17%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = '''3.33'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  +1.05%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = '''4.99'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  +4.67%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''12.76'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  +2.72%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('B31').Value = '''Hedera'
$ws.Range('B31').Style = 'Normal'
$ws.Range('C31').Value = '''https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('C31').Style = 'Normal'
$ws.Range('D31').Value = '''0.123'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  +10.17%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('B32').Value = '''Toncoin'
$ws.Range('B32').Style = 'Normal'
$ws.Range('C32').Value = '''https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('C32').Style = 'Normal'
$ws.Range('D32').Value = '''2.79'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  +4.70%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = '''7.36'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  -1.97%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = '''42.41'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  +6.30%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = '''  -0.90%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''1.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  +0.02%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = '''56.47'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  -1.83%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''0.0494'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  +1.31%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = '''0.0₃0760'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  +14.52%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''3.10'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  +33.66%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''0.147'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  +5.63%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = '''28.71'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  +31.52%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = '''0.996'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  -0.32%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = '''  +3.88%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '''2.98'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  -3.00%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = '''2.12'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  +6.15%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('B47').Value = '''NEARProtocol'
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').Value = '''https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').Value = '''4.45'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  -0.61%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('B48').Value = '''Monero'
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').Value = '''https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').Value = '''146.49'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  +1.40%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = '''  -3.09%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = '''0.311'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  -1.62%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''0.161'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  +17.06%  '
$ws.Range('E51').Style = 'Normal'
